$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Cells.Item(15, 8).Value = 1285.05
$ws.Cells.Item(15, 9).Value = 1285.05
$ws.Cells.Item(15, 11).Value = 3855.15
$ws.Cells.Item(15, 13).Value = -3686.15
# Row 34
$ws.Cells.Item(34, 8).Value = 451.33334
$ws.Cells.Item(34, 9).Value = 451.33334
$ws.Cells.Item(34, 11).Value = 451.33334
$ws.Cells.Item(34, 13).Value = -248.33334
# Row 36
$ws.Cells.Item(36, 8).Value = 451.33334
$ws.Cells.Item(36, 9).Value = 451.33334
$ws.Cells.Item(36, 11).Value = 451.33334
$ws.Cells.Item(36, 13).Value = 263.66666
# Row 40
$ws.Cells.Item(40, 8).Value = 1400
$ws.Cells.Item(40, 9).Value = 1500
$ws.Cells.Item(40, 10).Value = 1100
$ws.Cells.Item(40, 11).Value = 1500
$ws.Cells.Item(40, 12).Value = 1100
$ws.Cells.Item(40, 13).Value = -1325
$ws.Cells.Item(40, 14).Value = -1450
# Row 62
$ws.Cells.Item(62, 8).Value = 194449100
$ws.Cells.Item(62, 9).Value = 100007800
$ws.Cells.Item(62, 10).Value = 312500740
$ws.Cells.Item(62, 11).Value = 100007800
$ws.Cells.Item(62, 12).Value = 312500740
$ws.Cells.Item(62, 13).Value = -100007176
$ws.Cells.Item(62, 14).Value = -312501988
# Row 65
$ws.Cells.Item(65, 8).Value = 194449100
$ws.Cells.Item(65, 9).Value = 100007800
$ws.Cells.Item(65, 10).Value = 312500740
$ws.Cells.Item(65, 11).Value = 500039000
$ws.Cells.Item(65, 12).Value = 1562503700
$ws.Cells.Item(65, 13).Value = -500035880
$ws.Cells.Item(65, 14).Value = -1562509940
# Row 74
$ws.Cells.Item(74, 8).Value = 3839.375
$ws.Cells.Item(74, 9).Value = 2994.6155
$ws.Cells.Item(74, 10).Value = 7500
$ws.Cells.Item(74, 11).Value = 2994.6155
$ws.Cells.Item(74, 12).Value = 7500
$ws.Cells.Item(74, 13).Value = -2058.6155
$ws.Cells.Item(74, 14).Value = -9372
# Row 77
$ws.Cells.Item(77, 8).Value = 3839.375
$ws.Cells.Item(77, 9).Value = 2994.6155
$ws.Cells.Item(77, 10).Value = 7500
$ws.Cells.Item(77, 11).Value = 14973.0775
$ws.Cells.Item(77, 12).Value = 37500
$ws.Cells.Item(77, 13).Value = -10293.0775
$ws.Cells.Item(77, 14).Value = -46860
# Row 138
$ws.Cells.Item(138, 8).Value = 2431.8455
$ws.Cells.Item(138, 9).Value = 1673.3721
$ws.Cells.Item(138, 10).Value = 3035.8147
$ws.Cells.Item(138, 11).Value = 5020.1163
$ws.Cells.Item(138, 12).Value = 9107.444100000001
$ws.Cells.Item(138, 13).Value = 119.8837000000003
$ws.Cells.Item(138, 14).Value = -19387.4441

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 37
$ws.Cells.Item(37, 8).Value = 11795.462
$ws.Cells.Item(37, 9).Value = 6150
$ws.Cells.Item(37, 10).Value = 16634.428
$ws.Cells.Item(37, 11).Value = 6150
$ws.Cells.Item(37, 12).Value = 16634.428
$ws.Cells.Item(37, 13).Value = -5877
$ws.Cells.Item(37, 14).Value = -17180.428
# Row 45
$ws.Cells.Item(45, 8).Value = 716823.2
$ws.Cells.Item(45, 9).Value = 1001201.1
$ws.Cells.Item(45, 11).Value = 1001201.1
$ws.Cells.Item(45, 13).Value = -1000824.1
# Row 61
$ws.Cells.Item(61, 8).Value = 2716792
$ws.Cells.Item(61, 9).Value = 1191333.6
$ws.Cells.Item(61, 10).Value = 29412314
$ws.Cells.Item(61, 11).Value = 1191333.6
$ws.Cells.Item(61, 12).Value = 29412314
$ws.Cells.Item(61, 13).Value = -1191121.6
$ws.Cells.Item(61, 14).Value = -29412738

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 82
$ws.Cells.Item(82, 8).Value = 16027.357
$ws.Cells.Item(82, 9).Value = 2178.7144
$ws.Cells.Item(82, 10).Value = 29876
$ws.Cells.Item(82, 11).Value = 2178.7144
$ws.Cells.Item(82, 12).Value = 29876
$ws.Cells.Item(82, 13).Value = -1795.7144
$ws.Cells.Item(82, 14).Value = -30642
# Row 85
$ws.Cells.Item(85, 8).Value = 16027.357
$ws.Cells.Item(85, 9).Value = 2178.7144
$ws.Cells.Item(85, 10).Value = 29876
$ws.Cells.Item(85, 11).Value = 2178.7144
$ws.Cells.Item(85, 12).Value = 29876
$ws.Cells.Item(85, 13).Value = -852.7143999999998
$ws.Cells.Item(85, 14).Value = -32528
# Row 122
$ws.Cells.Item(122, 8).Value = 2821
$ws.Cells.Item(122, 9).Value = 1956
$ws.Cells.Item(122, 10).Value = 3253.5
$ws.Cells.Item(122, 11).Value = 5868
$ws.Cells.Item(122, 12).Value = 9760.5
$ws.Cells.Item(122, 13).Value = -3418
$ws.Cells.Item(122, 14).Value = -14660.5
# Row 134
$ws.Cells.Item(134, 8).Value = 15307397
$ws.Cells.Item(134, 9).Value = 20000992
$ws.Cells.Item(134, 10).Value = 3573409.8
$ws.Cells.Item(134, 11).Value = 60002976
$ws.Cells.Item(134, 12).Value = 10720229.4
$ws.Cells.Item(134, 13).Value = -60000441
$ws.Cells.Item(134, 14).Value = -10725299.4
# Row 136
$ws.Cells.Item(136, 8).Value = 2716792
$ws.Cells.Item(136, 9).Value = 1191333.6
$ws.Cells.Item(136, 10).Value = 29412314
$ws.Cells.Item(136, 11).Value = 3574000.8
$ws.Cells.Item(136, 12).Value = 88236942
$ws.Cells.Item(136, 13).Value = -3571450.8
$ws.Cells.Item(136, 14).Value = -88242042
# Row 139
$ws.Cells.Item(139, 8).Value = 51519
$ws.Cells.Item(139, 10).Value = 51377.777
$ws.Cells.Item(139, 12).Value = 51377.777
$ws.Cells.Item(139, 14).Value = -61657.777

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Cells.Item(16, 8).Value = 1858.7646
$ws.Cells.Item(16, 9).Value = 1574.9166
$ws.Cells.Item(16, 10).Value = 2540
$ws.Cells.Item(16, 11).Value = 1574.9166
$ws.Cells.Item(16, 12).Value = 2540
$ws.Cells.Item(16, 13).Value = -1287.9166
$ws.Cells.Item(16, 14).Value = -3114
# Row 109
$ws.Cells.Item(109, 8).Value = 18500
$ws.Cells.Item(109, 10).Value = 18500
$ws.Cells.Item(109, 12).Value = 18500
$ws.Cells.Item(109, 14).Value = -20580
# Row 113
$ws.Cells.Item(113, 8).Value = 1858.7646
$ws.Cells.Item(113, 9).Value = 1574.9166
$ws.Cells.Item(113, 10).Value = 2540
$ws.Cells.Item(113, 11).Value = 1574.9166
$ws.Cells.Item(113, 12).Value = 2540
$ws.Cells.Item(113, 13).Value = 595.0834
$ws.Cells.Item(113, 14).Value = -6880
# Row 132
$ws.Cells.Item(132, 8).Value = 1685.326
$ws.Cells.Item(132, 9).Value = 1426.4736
$ws.Cells.Item(132, 10).Value = 2914.875
$ws.Cells.Item(132, 11).Value = 4279.4208
$ws.Cells.Item(132, 12).Value = 8744.625
$ws.Cells.Item(132, 13).Value = -1749.4208
$ws.Cells.Item(132, 14).Value = -13804.625
# Row 133
$ws.Cells.Item(133, 8).Value = 39663
$ws.Cells.Item(133, 10).Value = 39663
$ws.Cells.Item(133, 12).Value = 39663
$ws.Cells.Item(133, 14).Value = -44723
# Row 135
$ws.Cells.Item(135, 8).Value = 43100
$ws.Cells.Item(135, 10).Value = 43100
$ws.Cells.Item(135, 12).Value = 43100
$ws.Cells.Item(135, 14).Value = -53240

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 23
$ws.Cells.Item(23, 8).Value = 5704.5
$ws.Cells.Item(23, 10).Value = 7880.077
$ws.Cells.Item(23, 12).Value = 23640.231
$ws.Cells.Item(23, 14).Value = -24110.231
# Row 109
$ws.Cells.Item(109, 8).Value = 3800.7222
$ws.Cells.Item(109, 9).Value = 883.3333
$ws.Cells.Item(109, 10).Value = 4165.396
$ws.Cells.Item(109, 11).Value = 2649.9999
$ws.Cells.Item(109, 12).Value = 12496.188
$ws.Cells.Item(109, 13).Value = -1609.9999
$ws.Cells.Item(109, 14).Value = -14576.188
# Row 131
$ws.Cells.Item(131, 8).Value = 7368467.5
$ws.Cells.Item(131, 10).Value = 960.2692
$ws.Cells.Item(131, 12).Value = 2880.8076
$ws.Cells.Item(131, 14).Value = -12960.8076

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 20
$ws.Cells.Item(20, 8).Value = 0
$ws.Cells.Item(20, 10).Value = 0
$ws.Cells.Item(20, 12).Value = 0
$ws.Range("N20").ClearContents()
# Row 61
$ws.Cells.Item(61, 8).Value = 1999
$ws.Cells.Item(61, 9).Value = 969.8570999999999
$ws.Cells.Item(61, 10).Value = 3800
$ws.Cells.Item(61, 11).Value = 969.8570999999999
$ws.Cells.Item(61, 12).Value = 3800
$ws.Cells.Item(61, 13).Value = -767.8570999999999
$ws.Cells.Item(61, 14).Value = -4204
# Row 113
$ws.Cells.Item(113, 8).Value = 1999
$ws.Cells.Item(113, 9).Value = 969.8570999999999
$ws.Cells.Item(113, 10).Value = 3800
$ws.Cells.Item(113, 11).Value = 969.8570999999999
$ws.Cells.Item(113, 12).Value = 3800
$ws.Cells.Item(113, 13).Value = 1200.1429
$ws.Cells.Item(113, 14).Value = -8140
# Row 122
$ws.Cells.Item(122, 8).Value = 23566426
$ws.Cells.Item(122, 9).Value = 15960088
$ws.Cells.Item(122, 10).Value = 66669000
$ws.Cells.Item(122, 11).Value = 47880264
$ws.Cells.Item(122, 12).Value = 200007000
$ws.Cells.Item(122, 13).Value = -47877814
$ws.Cells.Item(122, 14).Value = -200011900
# Row 140
$ws.Cells.Item(140, 8).Value = 58463.332
$ws.Cells.Item(140, 10).Value = 58463.332
$ws.Cells.Item(140, 12).Value = 58463.332
$ws.Cells.Item(140, 14).Value = -68823.33199999999

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 39
$ws.Cells.Item(39, 8).Value = 5735
$ws.Cells.Item(39, 9).Value = 2000
$ws.Cells.Item(39, 10).Value = 6980
$ws.Cells.Item(39, 11).Value = 2000
$ws.Cells.Item(39, 12).Value = 6980
$ws.Cells.Item(39, 13).Value = -1587
$ws.Cells.Item(39, 14).Value = -7806

Write-Output "All edits applied"